$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46082 -> 46083, i.e. 2026-03-01 -> 2026-03-02) for every data row (2..514).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 514) { $lastRow = 514 }

$ws.Range("C2:C$lastRow").Value = 46083
